$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version + Date, insert new "Jurisdiction" row ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$meta.Rows.Item(11).Insert()

# Match the formatting of the surrounding table rows
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Sheet "Elements": add II-1 constraint text to Subject.typeId row (row 5), column AJ (Constraint(s)) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}
"
